$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.35463285446167
$ws.Range("B1").Value = 1.533397436141968
$ws.Range("C1").Value = 1.297408819198608
$ws.Range("D1").Value = 1.354135155677795
$ws.Range("E1").Value = 0.9943990707397461
